$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4101661145687103
$ws.Range("B1").Value = 0.3308739364147186
$ws.Range("C1").Value = 0.2917959690093994
$ws.Range("D1").Value = 0.318186491727829
$ws.Range("E1").Value = 0.390886127948761
